# Lecture partielle de l'EDT M1 MIAGE.
# Update the three date rows (and their corresponding weekday labels)
# to a different week, keeping everything else unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 2023-03-11 (samedi) -> 2026-03-11 (mercredi)
$ws.Range("A2").Value = 46092.0
$ws.Range("B2").Value = "mercredi"

# Row 5: 2023-03-17 (vendredi) -> 2026-03-17 (mardi)
$ws.Range("A5").Value = 46098.0
$ws.Range("B5").Value = "mardi"

# Row 7: 2023-03-19 (dimanche) -> 2026-03-19 (jeudi)
$ws.Range("A7").Value = 46100.0
$ws.Range("B7").Value = "jeudi"
